$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell value for B8 with the new shared string "Ayudantia`nestructura"
$ws.Range("B8").Value = "Ayudantia`nestructura"

# Apply the same look as B4/B5 (green fill, centered, wrap text) to B8
$ws.Range("B8").WrapText = $True
$ws.Range("B8").Interior.Color = $ws.Range("B4").Interior.Color
$ws.Range("B8").HorizontalAlignment = $ws.Range("B4").HorizontalAlignment

# Update the selected cell to C6 (mirrors the selection change in the diff)
$ws.Range("C6").Select()
